# Chapter 6 - reveal the hidden (white-on-white) message by:
#  1. turning the hidden runs' font colour from white (FFFFFF) to red (FF0000)
#  2. replacing the plaintext with its substitution-cipher ciphertext

$d = $word.ActiveDocument

$replacements = @(
    @{
        Old = "The Columbian deal will be for 2 new venture wildcat wells, one each in the Llanos & Magdalena"
        New = "Fvr Gmxizfgmb qiyx kvpj ns ssp 2 zsj zczhhvc iwyhamh jijxg, brc qopl gz hui Jxoasq & Yothyxsae"
    },
    @{
        Old = "Basins.  These wells include a carry of thirty percent for the national oil company "
        New = "Zmgvrq.  Fvrwc isypq ubppsps n gydfl sd fvvvrk drvaqbg jmd hui lmhvslmz bmj ocztyzm "
    },
    @{
        Old = "and will test at least 3 K meters of vertical section.  In return, the client will be permitted "
        New = "nrb iwyp rqgg er xsnwr 3 W arxcdg bj tqfgmamz fiafwbr.  Gz frxsdb, glc ozvilf kvpj ns cipywgxcp "
    },
    @{
        Old = "to drill ten wells in the productive Putumayo province, earning a sixty % interest with a fifty "
        New = "hb hpuzy xcz krpje wa xfq desbgqgmtq Dhxsyols ndcimlos, repzwak y ewkxw % ubgipqgg agfv n jgrhl "
    },
    @{
        Old = "percent royalty rate, increasing to the standard eighty five percent royalty five years "
        New = "tcdqrrr dclejfm eerq, wagpqofmls hb xfq ggelpoeh cuuuxw rwii nqfpilf fbcyxhl jghs liydg "
    },
    @{
        Old = "after start of production in each well."
        New = "njrqf fxydh bj ndcqyafwbr gz sngf isyp."
    }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $range.Find.Font.Color = 16777215   # wdColorWhite (0xFFFFFF)
    $range.Find.Replacement.Font.Color = 255   # wdColorRed (0x0000FF -> BGR layout = 0xFF -> red)

    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.New, 2)
}
